# Updates the cryptocurrency price/volume table (columns D and E) for rows 2-51
# to match the latest scrape, per the GitHub Actions-driven commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.527.12"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "2.469.21"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  -0.43%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "314.51"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "91.88"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -2.25%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.38%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.514"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +3.21%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "32.51"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -3.17%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0792"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "2.850.41"
$ws.Range("E13").Value = "  -0.73%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.85"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.07%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "15.99"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +3.06%  "
$ws.Range("D16").Value = "2.480.68"
$ws.Range("E16").Value = "  -1.60%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.777"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").Value = "41.535.73"
$ws.Range("E18").Value = "  +0.16%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.50"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("D20").Value = "0.0₃0942"
$ws.Range("E20").Value = "  +0.47%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "70.85"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.68%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "11.08"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.46%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "237.89"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("E26").Value = "  +0.03%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "24.73"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("E28").Value = "  -1.40%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.69"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.76%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "35.38"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -4.52%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "155.62"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("E33").Value = "  -0.07%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0760"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.81%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "17.29"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -3.56%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.38"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.91%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.88"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -5.78%  "
$ws.Range("E38").Value = "  +2.97%  "
$ws.Range("E39").Value = "  +0.44%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.79"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -4.27%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "3.99"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -3.39%  "
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").Value = "1.945.81"
$ws.Range("E43").Value = "  -2.13%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0283"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.87%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "18.82"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -5.08%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.90"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -3.85%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "9.08"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("D48").Value = "2.709.16"
$ws.Range("E48").Value = "  -0.92%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "97.29"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.31%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "67.20"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -3.18%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "52.30"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +2.83%  "
